$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.270.22'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +1.50%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.940.30'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '492.59'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.10%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.30'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('E7').Value = '  -1.14%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.733'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.22%  '
$ws.Range('E10').Value = '  +3.31%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000352'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.80%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '43.26'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('E13').Value = '  -1.83%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.574.00'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.913.72'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.17%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.31'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -3.07%  '
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '19.86'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.68%  '
$ws.Range('E19').Value = '  +3.23%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '69.381.84'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.50%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '440.42'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.80%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '3.45'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.87%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '14.47'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.93%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '88.82'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '12.07'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +7.79%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.80'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +4.87%  '
$ws.Range('E27').Value = '  -4.66%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '37.11'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.71%  '
$ws.Range('E29').Value = '  -4.40%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '705.75'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -1.40%  '
$ws.Range('E31').Value = '  -0.23%  '
$ws.Range('E32').Value = '  -1.07%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '2.90'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.66%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.462'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +19.35%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0₃0912'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.96%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '6.09'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.36%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '61.51'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +4.16%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '40.84'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.20%  '
$ws.Range('E39').Value = '  +1.74%  '
$ws.Range('E40').Value = '  +0.01%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0489'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +1.66%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.91'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.34%  '
$ws.Range('E44').Value = '  -3.25%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.00'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +2.70%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.144'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +0.17%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.34'
$ws.Range('D47').ClearFormats()
$ws.Range('E48').Value = '  +3.87%  '
$ws.Range('E49').Value = '  +5.78%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '3.39'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.56%  '
$ws.Range('E51').Value = '  -0.80%  '
